$d = $word.ActiveDocument
$d.Content.Find.Execute("pusher sur cette branch ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "pusher sur cette branche", 2)
